$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowVals = @(-17.51839190784226, 2.112481309058413, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, 2.950583070020239, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, 2.262200721767418, 2.269167952927676, -17.51839190784226, 3.362689930566073, -17.51839190784226, 1.878892070424718, -17.51839190784226, 2.068347111316704, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, 1.814634562245732, -17.51839190784226, -17.51839190784226, -17.51839190784226, 2.992145314714005, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(2.6535777025783, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(7, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, 1.781980763268676, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(8, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(3.776909004622728, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(9, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, 1.63568696679679, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(10, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, 2.809578406509093, -17.51839190784226, 2.616717382920615, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(11, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(12, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, 2.433008074673277, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, 2.437866542345386, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(13, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, 1.356268317045268, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(14, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, 1.347598049625287, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(15, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, 2.562020981065354, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(16, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, 1.799835330766399, 2.092692342334562, -17.51839190784226, -17.51839190784226, -17.51839190784226, 0.992901338695752, 1.144085243397168, 1.139139384971257, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(17, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, -17.51839190784226, 0.8287666071746772, 0.6653064661143054, 1.18879531692606, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(18, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, -17.51839190784226, 1.590676091168865, -17.51839190784226, -17.51839190784226, -17.51839190784226, 1.608577178651664, 1.826296266155177, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(19, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, 0.8301504003066414, 1.48389845791878, -17.51839190784226, 3.279966766921729, -17.51839190784226, 2.123703867074184, 0.8637296883539985, -17.51839190784226, 4.321920794636185)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(20, 2 + $i).Value = $rowVals[$i] }

$rowVals = @(-17.51839190784226, 1.071058086419863, -17.51839190784226, 2.053095967466738, -17.51839190784226, 2.563113635289158, 2.364019791730794, -17.51839190784226, -17.51839190784226, -17.51839190784226)
for ($i = 0; $i -lt $rowVals.Length; $i++) { $ws.Cells.Item(21, 2 + $i).Value = $rowVals[$i] }

